# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates column C (td_sim_1) and column D (record_atd) values for rows 2-38,
# and the summary average in C39, per the corrected simulation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = @(99,26,196,308,119,194,46,142,61,31,127,53,217,226,14,24,49,38,11,56,87,35,118,12,106,33,65,57,80,10,2,179,25,20,97,37,70)
$newD = @(93,38,192.5,327,110,191,44.5,135.5,44.5,41,120.5,43,211,234.5,14,24.5,65,32.5,8,49,89.5,44.5,121,9,109.5,29.5,55,58,63,8.5,192,169.5,24,15,90.5,33,42.5)

for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newC[$i]
    $ws.Cells.Item($row, 4).Value = $newD[$i]
}

# Update the summary average of column C (row 39)
$ws.Cells.Item(39, 3).Value = 82.97297297297297
